# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column headers, styled like the existing headers.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

$ws.Range("AC1:AE1").Font.Bold = $true
$ws.Range("AC1:AE1").Borders.Weight = 2
$ws.Range("AC1:AE1").HorizontalAlignment = -4108
$ws.Range("AC1:AE1").VerticalAlignment = -4160

# Data rows (2-35) - every player on this roster shares the team's
# 1994 Houston season record: 66 wins, 49 losses, 0 ties.
$ws.Range("AC2:AC35").Value = 66
$ws.Range("AD2:AD35").Value = 49
$ws.Range("AE2:AE35").Value = 0
